$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update customer KH040 (row 41): name, phone number, and gender
$ws.Range("B41").Value = "Hồ Minh Hậu"
$ws.Range("C41").Value = "'0585576500"
$ws.Range("C41").ClearFormats()
$ws.Range("D41").Value = "Nam"
